$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 1 header: drop the bold font + thin border + center/top
# alignment that used to highlight the header row -> back to the
# default "Normal" style (matches styles.xml losing fontId 1 /
# borderId 1 / the second cellXfs entry).
# ------------------------------------------------------------------
$ws.Range("A1:O1").Style = "Normal"

# --- Row 1: header text: "gemini" column removed (cols shift left); "Unnamed: 0" / new trailing col cleared ---
$ws.Cells.Item(1,1).Value = ""
$ws.Cells.Item(1,1).Font.Bold = $false
$ws.Cells.Item(1,5).Value = "literals"
$ws.Cells.Item(1,6).Value = "method call 1"
$ws.Cells.Item(1,7).Value = "method call 2"
$ws.Cells.Item(1,8).Value = "method call 3"
$ws.Cells.Item(1,9).Value = "Method decleration"
$ws.Cells.Item(1,10).Value = "parameters"
$ws.Cells.Item(1,11).Value = "retrun statement"
$ws.Cells.Item(1,12).Value = "summary"
$ws.Cells.Item(1,13).Value = "var 1"
$ws.Cells.Item(1,14).Value = "var 2"
$ws.Cells.Item(1,15).Value = ""
$ws.Cells.Item(1,15).Font.Bold = $false

# --- Row 3: was "Dwells with fixations"; now "Revisit count" with its own data, shifted one col left ---
$ws.Cells.Item(3,1).Value = "Revisit count"
$ws.Cells.Item(3,2).Value = 8
$ws.Cells.Item(3,3).Value = 34
$ws.Cells.Item(3,4).Value = 4
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 2
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 3
$ws.Cells.Item(3,9).Value = 17
$ws.Cells.Item(3,10).Value = ""
$ws.Cells.Item(3,10).Font.Bold = $false
$ws.Cells.Item(3,11).Value = 16
$ws.Cells.Item(3,12).Value = 47
$ws.Cells.Item(3,13).Value = ""
$ws.Cells.Item(3,13).Font.Bold = $false
$ws.Cells.Item(3,14).Value = 0
$ws.Cells.Item(3,15).Value = ""
$ws.Cells.Item(3,15).Font.Bold = $false

# --- Row 4: was "Revisit count"; now "Fixation count" with its own data ---
$ws.Cells.Item(4,1).Value = "Fixation count"
$ws.Cells.Item(4,2).Value = 13
$ws.Cells.Item(4,3).Value = 92
$ws.Cells.Item(4,4).Value = 8
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 5
$ws.Cells.Item(4,9).Value = 23
$ws.Cells.Item(4,10).Value = ""
$ws.Cells.Item(4,10).Font.Bold = $false
$ws.Cells.Item(4,11).Value = 26
$ws.Cells.Item(4,12).Value = 134
$ws.Cells.Item(4,13).Value = ""
$ws.Cells.Item(4,13).Font.Bold = $false
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = ""
$ws.Cells.Item(4,15).Font.Bold = $false

# --- Row 5: was "Fixation count"; now "Dwell time (ms)" (new metric row) ---
$ws.Cells.Item(5,1).Value = "Dwell time (ms)"
$ws.Cells.Item(5,2).Value = 4070.91
$ws.Cells.Item(5,3).Value = 26436.41
$ws.Cells.Item(5,4).Value = 3011.58
$ws.Cells.Item(5,5).Value = 650.66
$ws.Cells.Item(5,6).Value = 567.3
$ws.Cells.Item(5,7).Value = 183.43
$ws.Cells.Item(5,8).Value = 2143.98
$ws.Cells.Item(5,9).Value = 8066.63
$ws.Cells.Item(5,10).Value = ""
$ws.Cells.Item(5,10).Font.Bold = $false
$ws.Cells.Item(5,11).Value = 8250.77
$ws.Cells.Item(5,12).Value = 34109.33
$ws.Cells.Item(5,13).Value = ""
$ws.Cells.Item(5,13).Font.Bold = $false
$ws.Cells.Item(5,14).Value = 175.14
$ws.Cells.Item(5,15).Value = ""
$ws.Cells.Item(5,15).Font.Bold = $false

# --- Row 6: brand new row: "Dwell time (%)" ---
$ws.Cells.Item(6,1).Value = "Dwell time (%)"
$ws.Cells.Item(6,2).Value = 3.46
$ws.Cells.Item(6,3).Value = 22.49
$ws.Cells.Item(6,4).Value = 2.56
$ws.Cells.Item(6,5).Value = 0.55
$ws.Cells.Item(6,6).Value = 0.48
$ws.Cells.Item(6,7).Value = 0.16
$ws.Cells.Item(6,8).Value = 1.82
$ws.Cells.Item(6,9).Value = 6.86
$ws.Cells.Item(6,10).Value = ""
$ws.Cells.Item(6,10).Font.Bold = $false
$ws.Cells.Item(6,11).Value = 7.02
$ws.Cells.Item(6,12).Value = 29.08
$ws.Cells.Item(6,13).Value = ""
$ws.Cells.Item(6,13).Font.Bold = $false
$ws.Cells.Item(6,14).Value = 0.15
$ws.Cells.Item(6,15).Value = ""
$ws.Cells.Item(6,15).Font.Bold = $false

# --- Row 7: was "Fixation duration (ms)" stats; now recomputed "Fixation duration (ms)" values ---
$ws.Cells.Item(7,1).Value = "Fixation duration (ms)"
$ws.Cells.Item(7,2).Value = 313.15
$ws.Cells.Item(7,3).Value = 287.35
$ws.Cells.Item(7,4).Value = 376.45
$ws.Cells.Item(7,5).Value = 325.33
$ws.Cells.Item(7,6).Value = 189.1
$ws.Cells.Item(7,7).Value = 183.43
$ws.Cells.Item(7,8).Value = 428.8
$ws.Cells.Item(7,9).Value = 350.72
$ws.Cells.Item(7,10).Value = ""
$ws.Cells.Item(7,10).Font.Bold = $false
$ws.Cells.Item(7,11).Value = 317.34
$ws.Cells.Item(7,12).Value = 254.55
$ws.Cells.Item(7,13).Value = ""
$ws.Cells.Item(7,13).Font.Bold = $false
$ws.Cells.Item(7,14).Value = 175.14
$ws.Cells.Item(7,15).Value = ""
$ws.Cells.Item(7,15).Font.Bold = $false

# --- Row 8: was "First fixation duration (ms)" (old row 7); values re-aligned to the new columns ---
$ws.Cells.Item(8,1).Value = "First fixation duration (ms)"
$ws.Cells.Item(8,2).Value = 166.89
$ws.Cells.Item(8,3).Value = 1001.11
$ws.Cells.Item(8,4).Value = 750.8099999999999
$ws.Cells.Item(8,5).Value = 166.86
$ws.Cells.Item(8,6).Value = 166.88
$ws.Cells.Item(8,7).Value = 183.43
$ws.Cells.Item(8,8).Value = 750.8099999999999
$ws.Cells.Item(8,9).Value = 116.77
$ws.Cells.Item(8,11).Value = 133.64
$ws.Cells.Item(8,12).Value = 350.37
$ws.Cells.Item(8,14).Value = 175.14
